$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.9968250993858209
$ws.Range("C3").Value = 0.9971944218932864
$ws.Range("D3").Value = 0.9963780486824455

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9973888076978041
$ws.Range("C4").Value = 0.9977317195496425
$ws.Range("D4").Value = 0.9977024659976855

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9977488042316179
$ws.Range("C5").Value = 0.9979096137536088
$ws.Range("D5").Value = 0.9979200096138564
